$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill A4:A65 with incrementing numbers 4..65
for ($i = 4; $i -le 65; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}

# Select A3:A65, with the active cell being A3
$ws.Range("A3:A65").Select()

# Scroll the window so that row 42 (A42) is the top-left visible cell
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1
